# Add bash entry about pushd/popd/dirs, sed replace and a case
# (appends rows 38-40 to the first worksheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$nl = [char]10

# --- Row 38: pushd/popd -------------------------------------------------
$ws.Range("A38").Value = 'pushd/popd'

# --- Row 39: sed ----------------------------------------------------------
$ws.Range("A39").Value = 'sed'
$ws.Range("C39").Value = 'echo $PATH | sed s/:/\\n/g'
$ws.Range("B39").Value = 'replace : character with \n'

# --- back to row 38 --------------------------------------------------------
$ws.Range("B38").Value = 'basic '
$ws.Range("C38").Value = 'pushd {path} //perform cd and push the path to stack' + $nl + 'popd {path} //pop the stack and cd to the poped path' + $nl + 'dirs //show your DIRectory Stack'

# --- Row 40: case study -----------------------------------------------------
$ws.Range("A40").Value = '案例'
$ws.Range("B40").Value = '找出bash在call哪個vim'
$ws.Range("C40").Value = 'for f in $(echo $PATH | sed s/:/\\n/g);' + $nl + '    do find $f -name "vim.exe";' + $nl + 'done'

# Wrap the long multi-line "CODE" column cells like the other such rows (style s="5")
$ws.Range("C38").WrapText = $true
$ws.Range("C40").WrapText = $true

# Row heights to fit the 3-line wrapped content (matches the 47.25pt used elsewhere)
$ws.Rows.Item(38).RowHeight = 47.25
$ws.Rows.Item(40).RowHeight = 47.25

# Move the active selection the same way Excel would after typing the new rows
$ws.Range("C41").Select()
